$d = $word.ActiveDocument

$replacements = @(
    @("398×8=3184", "169×9=1521"),
    @("493×5=2465", "314×9=2826"),
    @("300×3=900", "344×7=2408"),
    @("488×3=1464", "386×3=1158"),
    @("979×9=8811", "263×8=2104"),
    @("743×7=5201", "960×6=5760"),
    @("847×9=7623", "805×2=1610"),
    @("910×8=7280", "757×7=5299"),
    @("898×3=2694", "791×4=3164"),
    @("670×6=4020", "356×8=2848"),
    @("553×3=1659", "944×7=6608"),
    @("629×8=5032", "517×9=4653"),
    @("320×3=960", "149×9=1341"),
    @("935×5=4675", "556×8=4448"),
    @("643×2=1286", "983×6=5898"),
    @("175×6=1050", "713×4=2852"),
    @("590×8=4720", "994×4=3976"),
    @("566×9=5094", "757×9=6813"),
    @("712×3=2136", "295×4=1180"),
    @("314×2=628", "359×3=1077"),
    @("950×7=6650", "246×2=492"),
    @("683×3=2049", "410×5=2050"),
    @("700×2=1400", "698×5=3490"),
    @("951×4=3804", "295×7=2065"),
    @("281×8=2248", "190×5=950")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
